# feat: add 2022-Q1 data
#
# Plan:
#  - duplicate the existing "总计" sheet (so the copy inherits the same
#    sheet-level setup: page margins, outline props, header style, etc.),
#  - rename the original to "2022-Q1" and turn it into a fund-holdings
#    sheet with the same column layout as "2021-Q4",
#  - rename the duplicate back to "总计" and rebuild it with the new
#    2022-Q1 summary row stacked on top of the pre-existing 2021-Q4 row.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate "总计" (placed right after itself), then rename the two
#    copies so the tab order becomes 2021-Q4, 2022-Q1, 总计.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Copy($null, $q1)
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Item("总计 (2)")
$total.Name = "总计"

# ------------------------------------------------------------------
# 2. Build the "2022-Q1" sheet - same header/column layout as "2021-Q4".
# ------------------------------------------------------------------

# headers (row 1) already carry the old "总计" header style (bold + thin
# border); extend that same style across the new E1:H1 header cells
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# row-label column A already carries its style on A2; extend it to A3
$q1.Range("A2").Copy()
$q1.Range("A3").PasteSpecial(-4122)

# row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'006143"
$q1.Range("C2").Value = "恒生前海中证质量成长低波动指数A"
$q1.Range("D2").Value = "'0.06"
$q1.Range("E2").Value = "'94.34"
$q1.Range("F2").Value = "'2.67"
$q1.Range("G2").Value = "'0.0016"
$q1.Range("H2").Value = 6

# row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'006144"
$q1.Range("C3").Value = "恒生前海中证质量成长低波动指数C"
$q1.Range("D3").Value = "'0.01"
$q1.Range("E3").Value = "'94.34"
$q1.Range("F3").Value = "'2.67"
$q1.Range("G3").Value = "'0.0003"
$q1.Range("H3").Value = 6

# ------------------------------------------------------------------
# 3. Rebuild the new "总计" sheet - 2022-Q1 stacked on top of the
#    pre-existing 2021-Q4 summary row.
# ------------------------------------------------------------------

# row-label column A style, extended down to row 3
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

# row 2 - 2022-Q1 (new)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

# row 3 - 2021-Q4 (previously row 2 of the old "总计" sheet)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0

$wb.Worksheets.Item("2021-Q4").Select()
